# testdata.xlsx edit
#
# DATA sheet ("DATA", the active/visible tab): the "execute" flag (column B)
# for the amazonTest / chrome / 88.0.4324.96 row (row 8) flips from "no" to
# "yes" - mirrors the Java-side RemoveIf optimisation now picking this case
# back up for execution.
#
# The selections left behind in both sheets' sheetViews are updated to match
# where the author was navigating: row 7:8 highlighted on DATA (the row that
# was just edited), and cell A4 left selected on RUNMANAGER.

$wb = $excel.ActiveWorkbook

# RUNMANAGER sheet - just leave the selection parked on A4.
$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsRunManager.Range("A4").Select()

# DATA sheet - the actual edit, then leave A7:H8 selected.
$wsData = $wb.Worksheets.Item("DATA")
$wsData.Activate()
$wsData.Range("B8").Value = "yes"
$wsData.Range("A7:H8").Select()
